$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3737415969371796
$ws.Range("B1").Value = 1.21735942363739
$ws.Range("C1").Value = 4.909544467926025
$ws.Range("D1").Value = 1.512856125831604
$ws.Range("E1").Value = 0.7582865357398987
